$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.997.03"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'  +0.06%  "
$ws.Range("E2").ClearFormats()

$ws.Range("D3").Value = "'1.885.26"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'  -1.68%  "
$ws.Range("E3").ClearFormats()

$ws.Range("E4").Value = "'  +0.18%  "
$ws.Range("E4").ClearFormats()

$ws.Range("D5").Value = "'331.10"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'  -2.27%  "
$ws.Range("E5").ClearFormats()

$ws.Range("D6").Value = "'1.002"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "'  +0.12%  "
$ws.Range("E6").ClearFormats()

$ws.Range("D7").Value = "'0.4600"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "'  -3.11%  "
$ws.Range("E7").ClearFormats()

$ws.Range("D8").Value = "'0.4061"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "'  -0.07%  "
$ws.Range("E8").ClearFormats()

$ws.Range("E9").Value = "'  -0.91%  "
$ws.Range("E9").ClearFormats()

$ws.Range("D10").Value = "'0.07981"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "'  -2.46%  "
$ws.Range("E10").ClearFormats()

$ws.Range("D11").Value = "'0.9902"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "'  -4.00%  "
$ws.Range("E11").ClearFormats()

$ws.Range("D12").Value = "'21.66"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "'  -3.67%  "
$ws.Range("E12").ClearFormats()

$ws.Range("D13").Value = "'1.886.76"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "'  -0.64%  "
$ws.Range("E13").ClearFormats()

$ws.Range("D14").Value = "'5.902"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "'  -3.16%  "
$ws.Range("E14").ClearFormats()

$ws.Range("D15").Value = "'7.063"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "'  -4.49%  "
$ws.Range("E15").ClearFormats()

$ws.Range("D16").Value = "'1.003"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "'  +0.20%  "
$ws.Range("E16").ClearFormats()

$ws.Range("D17").Value = "'88.44"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "'  -3.70%  "
$ws.Range("E17").ClearFormats()

$ws.Range("D18").Value = "'0.00001031"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "'  -2.28%  "
$ws.Range("E18").ClearFormats()

$ws.Range("D19").Value = "'0.06560"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "'  -1.10%  "
$ws.Range("E19").ClearFormats()

$ws.Range("E20").Value = "'  -2.65%  "
$ws.Range("E20").ClearFormats()

$ws.Range("E21").Value = "'  -0.02%  "
$ws.Range("E21").ClearFormats()

$ws.Range("D22").Value = "'29.025.92"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "'  +0.08%  "
$ws.Range("E22").ClearFormats()

$ws.Range("D23").Value = "'5.420"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "'  -3.02%  "
$ws.Range("E23").ClearFormats()

$ws.Range("D24").Value = "'11.42"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "'  +2.03%  "
$ws.Range("E24").ClearFormats()

$ws.Range("D25").Value = "'2.209"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "'  -2.73%  "
$ws.Range("E25").ClearFormats()

$ws.Range("D26").Value = "'2.098.75"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "'  -1.35%  "
$ws.Range("E26").ClearFormats()

$ws.Range("D27").Value = "'156.96"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "'  -2.44%  "
$ws.Range("E27").ClearFormats()

$ws.Range("D28").Value = "'19.59"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "'  -2.46%  "
$ws.Range("E28").ClearFormats()

$ws.Range("D29").Value = "'2.091"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "'  -4.60%  "
$ws.Range("E29").ClearFormats()

$ws.Range("D30").Value = "'5.488"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "'  -1.00%  "
$ws.Range("E30").ClearFormats()

$ws.Range("D31").Value = "'117.49"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "'  -3.02%  "
$ws.Range("E31").ClearFormats()

$ws.Range("D32").Value = "'1.013"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "'  +0.04%  "
$ws.Range("E32").ClearFormats()

$ws.Range("D33").Value = "'0.09330"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "'  -2.64%  "
$ws.Range("E33").ClearFormats()

$ws.Range("D34").Value = "'3.604"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "'  -1.04%  "
$ws.Range("E34").ClearFormats()

$ws.Range("D35").Value = "'1.404"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "'  -2.61%  "
$ws.Range("E35").ClearFormats()

$ws.Range("D36").Value = "'5.275"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "'  -3.01%  "
$ws.Range("E36").ClearFormats()

$ws.Range("D37").Value = "'0.06061"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "'  -2.40%  "
$ws.Range("E37").ClearFormats()

$ws.Range("D38").Value = "'0.02219"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "'  -2.86%  "
$ws.Range("E38").ClearFormats()

$ws.Range("D39").Value = "'8.284"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "'  -4.88%  "
$ws.Range("E39").ClearFormats()

$ws.Range("D40").Value = "'1.175"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "'  -2.25%  "
$ws.Range("E40").ClearFormats()

$ws.Range("E41").Value = "'  +0.12%  "
$ws.Range("E41").ClearFormats()

$ws.Range("D42").Value = "'0.5782"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "'  -4.44%  "
$ws.Range("E42").ClearFormats()

$ws.Range("B43").Value = "'Algorand"
$ws.Range("B43").ClearFormats()
$ws.Range("C43").Value = "'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("C43").ClearFormats()
$ws.Range("D43").Value = "'0.1824"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "'  -4.13%  "
$ws.Range("E43").ClearFormats()

$ws.Range("B44").Value = "'Aptos"
$ws.Range("B44").ClearFormats()
$ws.Range("C44").Value = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("C44").ClearFormats()
$ws.Range("D44").Value = "'10.13"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "'  -4.57%  "
$ws.Range("E44").ClearFormats()

$ws.Range("D45").Value = "'1.264"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "'  -1.34%  "
$ws.Range("E45").ClearFormats()

$ws.Range("D46").Value = "'0.07442"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "'  +2.63%  "
$ws.Range("E46").ClearFormats()

$ws.Range("B47").Value = "'EnergySwap"
$ws.Range("B47").ClearFormats()
$ws.Range("C47").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("C47").ClearFormats()
$ws.Range("D47").Value = "'12.07"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "'  -2.39%  "
$ws.Range("E47").ClearFormats()

$ws.Range("B48").Value = "'RenderToken"
$ws.Range("B48").ClearFormats()
$ws.Range("C48").Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("C48").ClearFormats()
$ws.Range("D48").Value = "'2.258"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "'  +4.72%  "
$ws.Range("E48").ClearFormats()

$ws.Range("D49").Value = "'0.5450"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "'  -3.40%  "
$ws.Range("E49").ClearFormats()

$ws.Range("D50").Value = "'1.897"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "'  -4.35%  "
$ws.Range("E50").ClearFormats()

$ws.Range("D51").Value = "'45.76"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "'  +11.95%  "
$ws.Range("E51").ClearFormats()
